# Apply the recorded edit to the workbook.
$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("RawPoints")
$wsPaths = $wb.Worksheets.Item("Paths")

# ---------------------------------------------------------------------------
# RawPoints sheet: update the underlying measurement values / formulas.
# All of the dependent H-column (and Paths sheet) values recompute from
# these automatically.
# ---------------------------------------------------------------------------

# C8 stays a static value, but its number changes.
$wsRaw.Range("C8").Value = 324

# C9 used to be a formula (=24*12); it becomes a plain static value.
$wsRaw.Range("C9").Value = 0

# C10 / C11 keep being formulas, but the sign of the offset flips.
$wsRaw.Range("C10").Formula = "=C8-30"
$wsRaw.Range("C11").Formula = "=C9+30"

# C18 / C19 / C21 / C25 / C26 / C27 / C28 get new static values, and each
# gains a companion "mirror" formula in column D (=324-C#).
$wsRaw.Range("C18").Value = 222
$wsRaw.Range("D18").Formula = "=324-C18"

$wsRaw.Range("C19").Value = 174
$wsRaw.Range("D19").Formula = "=324-C19"

$wsRaw.Range("C21").Value = 238.75
$wsRaw.Range("D21").Formula = "=324-C21"

# C22 references C8 instead of C9 now.
$wsRaw.Range("C22").Formula = "=C8-C21"

$wsRaw.Range("C25").Value = 252.43
$wsRaw.Range("D25").Formula = "=324-C25"

$wsRaw.Range("C26").Value = 217.75
$wsRaw.Range("D26").Formula = "=324-C26"

$wsRaw.Range("C27").Value = 71.57
$wsRaw.Range("D27").Formula = "=324-C27"

$wsRaw.Range("C28").Value = 107.57
$wsRaw.Range("D28").Formula = "=324-C28"

# ---------------------------------------------------------------------------
# Paths sheet: highlight row 4 in yellow, and insert a new "Move8" row.
# ---------------------------------------------------------------------------

$wsPaths.Range("A4:R4").Interior.Color = 65535

$wsPaths.Rows("31:31").Insert()

$wsPaths.Range("A31").Value = "Move8"
$wsPaths.Range("B31").Value = "Move eight feet"
$wsPaths.Range("C31").Value = 2
$wsPaths.Range("D31").Value = 0
$wsPaths.Range("E31").Value = 0
$wsPaths.Range("F31").Value = 0
$wsPaths.Range("G31").Value = 96
$wsPaths.Range("H31").Value = 0
$wsPaths.Range("I31").Value = 0

# ---------------------------------------------------------------------------
# Selection / active-sheet state.
# ---------------------------------------------------------------------------

$wsRaw.Activate()
$wsRaw.Range("C33").Select()

$wsPaths.Activate()
$wsPaths.Range("G32").Select()
